$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'" + '29.503.41'
$ws.Range("E2").Value = '  -2.87%  '
$ws.Range("D3").Value = "'" + '1.994.89'
$ws.Range("E3").Value = '  -4.82%  '
$ws.Range("D4").Value = "'" + '1.015'
$ws.Range("E4").Value = '  +1.19%  '
$ws.Range("D5").Value = "'" + '329.48'
$ws.Range("E5").Value = '  -3.95%  '
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("D7").Value = "'" + '0.5022'
$ws.Range("E7").Value = '  -4.54%  '
$ws.Range("D8").Value = "'" + '0.4238'
$ws.Range("E8").Value = '  -4.34%  '
$ws.Range("D9").Value = "'" + '53.72'
$ws.Range("E9").Value = '  -1.68%  '
$ws.Range("D10").Value = "'" + '0.08950'
$ws.Range("E10").Value = '  -4.28%  '
$ws.Range("D11").Value = "'" + '1.113'
$ws.Range("E11").Value = '  -4.88%  '
$ws.Range("D12").Value = "'" + '23.21'
$ws.Range("E12").Value = '  -6.26%  '
$ws.Range("D13").Value = "'" + '2.001.18'
$ws.Range("E13").Value = '  -6.01%  '
$ws.Range("D14").Value = "'" + '7.978'
$ws.Range("E14").Value = '  -6.81%  '
$ws.Range("D15").Value = "'" + '6.461'
$ws.Range("E15").Value = '  -6.63%  '
$ws.Range("E16").Value = '  +1.13%  '
$ws.Range("D17").Value = "'" + '94.41'
$ws.Range("E17").Value = '  -6.98%  '
$ws.Range("D18").Value = "'" + '0.00001113'
$ws.Range("E18").Value = '  -4.04%  '
$ws.Range("D19").Value = "'" + '0.06755'
$ws.Range("E19").Value = '  +0.95%  '
$ws.Range("D20").Value = "'" + '19.42'
$ws.Range("E20").Value = '  -8.49%  '
$ws.Range("E21").Value = '  +1.06%  '
$ws.Range("D22").Value = "'" + '5.940'
$ws.Range("E22").Value = '  -6.30%  '
$ws.Range("D23").Value = "'" + '29.537.43'
$ws.Range("E23").Value = '  -2.90%  '
$ws.Range("D24").Value = "'" + '12.09'
$ws.Range("E24").Value = '  -3.83%  '
$ws.Range("D25").Value = "'" + '2.316'
$ws.Range("E25").Value = '  +0.30%  '
$ws.Range("D26").Value = "'" + '20.80'
$ws.Range("E26").Value = '  -4.84%  '
$ws.Range("D27").Value = "'" + '157.10'
$ws.Range("E27").Value = '  -3.73%  '
$ws.Range("D28").Value = "'" + '6.307'
$ws.Range("E28").Value = '  -7.27%  '
$ws.Range("D29").Value = "'" + '2.305'
$ws.Range("E29").Value = '  -8.37%  '
$ws.Range("E30").Value = '  -4.53%  '
$ws.Range("D31").Value = "'" + '1.060'
$ws.Range("E31").Value = '  -7.03%  '
$ws.Range("D32").Value = "'" + '0.09941'
$ws.Range("E32").Value = '  -5.36%  '
$ws.Range("D33").Value = "'" + '1.548'
$ws.Range("E33").Value = '  -6.41%  '
$ws.Range("D34").Value = "'" + '5.834'
$ws.Range("E34").Value = '  -6.98%  '
$ws.Range("E35").Value = '  -1.99%  '
$ws.Range("D36").Value = "'" + '0.02463'
$ws.Range("E36").Value = '  -6.61%  '
$ws.Range("D37").Value = "'" + '9.247'
$ws.Range("E37").Value = '  -9.00%  '
$ws.Range("D38").Value = "'" + '0.06396'
$ws.Range("E38").Value = '  -5.94%  '
$ws.Range("D39").Value = "'" + '1.297'
$ws.Range("E39").Value = '  -3.60%  '
$ws.Range("D40").Value = "'" + '0.6537'
$ws.Range("E40").Value = '  -6.80%  '
$ws.Range("D41").Value = "'" + '11.62'
$ws.Range("E41").Value = '  -7.81%  '
$ws.Range("D42").Value = "'" + '0.2040'
$ws.Range("E42").Value = '  -8.19%  '
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("D44").Value = "'" + '0.6332'
$ws.Range("E44").Value = '  -7.85%  '
$ws.Range("B45").Value = 'NEARProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D45").Value = "'" + '2.211'
$ws.Range("E45").Value = '  -5.68%  '
$ws.Range("B46").Value = 'EnergySwap'
$ws.Range("C46").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D46").Value = "'" + '13.44'
$ws.Range("E46").Value = '  -6.53%  '
$ws.Range("D47").Value = "'" + '1.313'
$ws.Range("E47").Value = '  -5.20%  '
$ws.Range("D48").Value = "'" + '3.502'
$ws.Range("E48").Value = '  -3.68%  '
$ws.Range("D49").Value = "'" + '0.00000000337'
$ws.Range("E49").Value = '  -2.01%  '
$ws.Range("D50").Value = "'" + '0.06958'
$ws.Range("E50").Value = '  -4.09%  '
$ws.Range("D51").Value = "'" + '1.132'
